$d = $word.ActiveDocument

# The last paragraph in the document body (before the final sectPr) is an
# empty paragraph styled in Times New Roman. We add right alignment and
# insert the text "Wojciech Dreslerski" as a run using that paragraph's
# formatting.
$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCount)

# Right-align the paragraph
$lastPara.Range.ParagraphFormat.Alignment = 2  # wdAlignParagraphRight

# Set the text of the paragraph (this will appear before the paragraph mark)
$lastPara.Range.Text = "Wojciech Dreslerski"

# Ensure the run uses Times New Roman, matching the existing formatting
# (ascii, hAnsi (NameOther) and cs (NameBi) font faces)
$lastPara.Range.Font.NameAscii = "Times New Roman"
$lastPara.Range.Font.NameOther = "Times New Roman"
$lastPara.Range.Font.NameBi = "Times New Roman"

$d.Save()
